$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the table (header row stays the same).
$data = @(
    @("Broadmeadows", "Craigieburn Line train", "1.25pm - 1.59pm  9/02/2012", "Case caught train from Broadmeadows Railway Station to Glenroy Railway Station", "old"),
    @("Broadmeadows", "Craigieburn Line train", "1.25pm - 1.59pm  9/02/2021", "Case caught train from Broadmeadows Railway Station to Glenroy Railway Station", "new"),
    @("Broadmeadows", "Ferguson Plarre Bakehouses - Broadmeadows, 1099-1169 Pascoe Vale Road", "12:30pm - 12:45pm 9/2/2021", "Case attended venue", "new"),
    @("Broadmeadows", "Woolworths Broadmeadows Central, Pascoe Vale Road", "12.15pm - 12:30 pm 9/2/2021", "Case attended venue", "new"),
    @("Keysborough", "Aces Sporting Club (Driving Range)  Cnr Springvale Rd and Hutton Rd  Keysborough VIC 3173", "10:00pm - 11:15pm  30/1/2021", "Case attended venue", "old"),
    @("Noble Park", "Club Noble  46/56 Moodemere St  Noble Park VIC 3174", "2:36pm -3:30pm  30/01/2021", "Case attended venue", "old"),
    @("Pascoe Vale", "Elite Swimming Pascoe Vale, 8 Attercliffe Avenue", "5pm - 6pm 8/2/2021", "Case attended venue", "new"),
    @("Pascoe Vale", "Oak Park Sports and Aquatic Centre, 563a Pascoe Vale Road", "4pm - 7.30pm 10/2/2021", "Case attended venue", "new")
)

# Remove the old data rows (2 through the last used row) before writing the
# replacement rows, so no stray rows remain if the old table was longer.
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -ge 2) {
    $ws.Range("A2:E" + $lastRow).ClearContents()
}

$rowIndex = 2
foreach ($record in $data) {
    $ws.Range("A" + $rowIndex).Value = $record[0]
    $ws.Range("B" + $rowIndex).Value = $record[1]
    $ws.Range("C" + $rowIndex).Value = $record[2]
    $ws.Range("D" + $rowIndex).Value = $record[3]
    $ws.Range("E" + $rowIndex).Value = $record[4]
    $rowIndex++
}

# The new table is shorter than the old one, so remove the now-unused
# trailing rows entirely (shifting cells up rather than leaving blanks).
$newLastRow = $rowIndex - 1
if ($lastRow -gt $newLastRow) {
    $ws.Range("A" + ($newLastRow + 1) + ":E" + $lastRow).EntireRow.Delete()
}

$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(5).EntireColumn.AutoFit()

$ws.Range("B8").Select()
